$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("normDataSource")

# Row 9: change "X" -> "XX" for all populated cells except F9 (which is blank)
$ws.Range("C9").Value = "XX"
$ws.Range("D9").Value = "XX"
$ws.Range("E9").Value = "XX"
$ws.Range("G9").Value = "XX"
$ws.Range("H9").Value = "XX"
$ws.Range("I9").Value = "XX"

# Row 10: change "X" -> "XX" for the populated cells (C10, G10)
$ws.Range("C10").Value = "XX"
$ws.Range("G10").Value = "XX"

# Update the active selection to C11
$ws.Activate()
$ws.Range("C11").Select()
